# "update for general search grid"
# Adds three new columns (L, M, N) to the TestData sheet to support a
# general search grid: GridTitle / SearchCriteria / SearchType headers
# in row 1, and PX / <blank> / first values in the data row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header) new cells ---
$ws.Range("L1").Value2 = "GridTitle"
$ws.Range("M1").Value2 = "SearchCriteria"
$ws.Range("N1").Value2 = "SearchType"

# Match the bold/centered header formatting used by the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# --- Row 4 (data) new cells ---
$ws.Range("L4").Value2 = "PX"
# M4 stays blank (matches the diff's empty <c r="M4" .../> cell).
$ws.Range("N4").Value2 = "first"

# Match the wrapped-text body formatting used by the rest of row 4.
$ws.Range("K4").Copy()
$ws.Range("L4:N4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Selection moves to the new last cell, N4.
$ws.Range("N4").Select()
